$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5557.5
$ws.Range("I62").Value = 5669.2
$ws.Range("K62").Value = 5669.2
$ws.Range("M62").Value = -5045.2
$ws.Range("H65").Value = 5557.5
$ws.Range("I65").Value = 5669.2
$ws.Range("K65").Value = 28346
$ws.Range("M65").Value = -25226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 30977.234
$ws.Range("I2").Value = 43169.207
$ws.Range("J2").Value = 1716.5
$ws.Range("K2").Value = 43169.207
$ws.Range("L2").Value = 1716.5
$ws.Range("M2").Value = -43056.207
$ws.Range("N2").Value = -1942.5
$ws.Range("H32").Value = 8835.853999999999
$ws.Range("I32").Value = 5329.7646
$ws.Range("J32").Value = 16286.292
$ws.Range("K32").Value = 5329.7646
$ws.Range("L32").Value = 16286.292
$ws.Range("M32").Value = -5042.7646
$ws.Range("N32").Value = -16860.292
$ws.Range("H61").Value = 2753.7144
$ws.Range("J61").Value = 4497.1665
$ws.Range("L61").Value = 4497.1665
$ws.Range("N61").Value = -4921.1665
$ws.Range("H74").Value = 50583.027
$ws.Range("I74").Value = 28490.883
$ws.Range("K74").Value = 28490.883
$ws.Range("M74").Value = -27616.883
$ws.Range("H77").Value = 50583.027
$ws.Range("I77").Value = 28490.883
$ws.Range("K77").Value = 142454.415
$ws.Range("M77").Value = -138086.415
$ws.Range("H97").Value = 4915.2383
$ws.Range("I97").Value = 5383.3613
$ws.Range("K97").Value = 5383.3613
$ws.Range("M97").Value = -4887.3613
$ws.Range("H116").Value = 30977.234
$ws.Range("I116").Value = 43169.207
$ws.Range("J116").Value = 1716.5
$ws.Range("K116").Value = 43169.207
$ws.Range("L116").Value = 1716.5
$ws.Range("M116").Value = -40875.207
$ws.Range("N116").Value = -6304.5
$ws.Range("H136").Value = 2753.7144
$ws.Range("J136").Value = 4497.1665
$ws.Range("L136").Value = 13491.4995
$ws.Range("N136").Value = -18591.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 30977.234
$ws.Range("I3").Value = 43169.207
$ws.Range("J3").Value = 1716.5
$ws.Range("K3").Value = 43169.207
$ws.Range("L3").Value = 1716.5
$ws.Range("M3").Value = -43055.207
$ws.Range("N3").Value = -1944.5
$ws.Range("H20").Value = 27781976
$ws.Range("I20").Value = 37039856
$ws.Range("K20").Value = 37039856
$ws.Range("M20").Value = -37039609
$ws.Range("H80").Value = 472.7742
$ws.Range("J80").Value = 498
$ws.Range("L80").Value = 498
$ws.Range("N80").Value = -2494
$ws.Range("H83").Value = 472.7742
$ws.Range("J83").Value = 498
$ws.Range("L83").Value = 2490
$ws.Range("N83").Value = -12474
$ws.Range("H86").Value = 12414.818
$ws.Range("J86").Value = 12358.091
$ws.Range("L86").Value = 12358.091
$ws.Range("N86").Value = -14604.091
$ws.Range("H89").Value = 12414.818
$ws.Range("J89").Value = 12358.091
$ws.Range("L89").Value = 61790.455
$ws.Range("N89").Value = -73022.455
$ws.Range("H107").Value = 2053.4285
$ws.Range("I107").Value = 1340.7273
$ws.Range("J107").Value = 4666.6665
$ws.Range("K107").Value = 1340.7273
$ws.Range("L107").Value = 4666.6665
$ws.Range("M107").Value = 579.2727
$ws.Range("N107").Value = -8506.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4225.3335
$ws.Range("I58").Value = 4358.8423
$ws.Range("J58").Value = 3718
$ws.Range("K58").Value = 4358.8423
$ws.Range("L58").Value = 3718
$ws.Range("M58").Value = -4155.8423
$ws.Range("N58").Value = -4124
$ws.Range("H105").Value = 5392.375
$ws.Range("I105").Value = 5313.3076
$ws.Range("J105").Value = 5735
$ws.Range("K105").Value = 5313.3076
$ws.Range("L105").Value = 5735
$ws.Range("M105").Value = -3566.3076
$ws.Range("N105").Value = -9229
$ws.Range("H122").Value = 2203.8572
$ws.Range("I122").Value = 2019.7333
$ws.Range("J122").Value = 2664.1667
$ws.Range("K122").Value = 6059.199900000001
$ws.Range("L122").Value = 7992.500100000001
$ws.Range("M122").Value = -3609.199900000001
$ws.Range("N122").Value = -12892.5001
$ws.Range("H134").Value = 22037.844
$ws.Range("I134").Value = 27373.885
$ws.Range("K134").Value = 82121.655
$ws.Range("M134").Value = -79586.655
$ws.Range("H136").Value = 4225.3335
$ws.Range("I136").Value = 4358.8423
$ws.Range("J136").Value = 3718
$ws.Range("K136").Value = 13076.5269
$ws.Range("L136").Value = 11154
$ws.Range("M136").Value = -10526.5269
$ws.Range("N136").Value = -16254

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1907
$ws.Range("I22").Value = 1749.5
$ws.Range("K22").Value = 5248.5
$ws.Range("M22").Value = -5079.5
$ws.Range("H27").Value = 1907
$ws.Range("I27").Value = 1749.5
$ws.Range("K27").Value = 5248.5
$ws.Range("M27").Value = -5146.5
$ws.Range("H56").Value = 16671847
$ws.Range("I56").Value = 16671847
$ws.Range("K56").Value = 16671847
$ws.Range("M56").Value = -16671317
$ws.Range("H117").Value = 1331.9
$ws.Range("I117").Value = 808
$ws.Range("J117").Value = 1424.3529
$ws.Range("K117").Value = 2424
$ws.Range("L117").Value = 4273.0587
$ws.Range("M117").Value = 1018
$ws.Range("N117").Value = -11157.0587
$ws.Range("H131").Value = 12628873
$ws.Range("I131").Value = 5556079
$ws.Range("J131").Value = 18522868
$ws.Range("K131").Value = 16668237
$ws.Range("L131").Value = 55568604
$ws.Range("M131").Value = -16663197
$ws.Range("N131").Value = -55578684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10589.25
$ws.Range("I70").Value = 11376.429
$ws.Range("K70").Value = 11376.429
$ws.Range("M70").Value = -11106.429
$ws.Range("H73").Value = 10589.25
$ws.Range("I73").Value = 11376.429
$ws.Range("K73").Value = 11376.429
$ws.Range("M73").Value = -10440.429
$ws.Range("H102").Value = 48624.727
$ws.Range("I102").Value = 2191.5454
$ws.Range("J102").Value = 95057.91
$ws.Range("K102").Value = 2191.5454
$ws.Range("L102").Value = 95057.91
$ws.Range("M102").Value = -569.5454
$ws.Range("N102").Value = -98301.91
$ws.Range("H132").Value = 2448.0667
$ws.Range("I132").Value = 2337.3096
$ws.Range("J132").Value = 3998.6667
$ws.Range("K132").Value = 7011.9288
$ws.Range("L132").Value = 11996.0001
$ws.Range("M132").Value = -4481.9288
$ws.Range("N132").Value = -17056.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6673.647
$ws.Range("I7").Value = 5332.9165
$ws.Range("K7").Value = 5332.9165
$ws.Range("M7").Value = -5220.9165
$ws.Range("H40").Value = 5490.4614
$ws.Range("I40").Value = 4108.3335
$ws.Range("J40").Value = 8600.25
$ws.Range("K40").Value = 4108.3335
$ws.Range("L40").Value = 8600.25
$ws.Range("M40").Value = -3972.3335
$ws.Range("N40").Value = -8872.25
$ws.Range("H46").Value = 7849.846
$ws.Range("J46").Value = 9118.25
$ws.Range("L46").Value = 9118.25
$ws.Range("N46").Value = -9494.25
$ws.Range("H82").Value = 66669500
$ws.Range("I82").Value = 111113700
$ws.Range("J82").Value = 3212.5
$ws.Range("K82").Value = 111113700
$ws.Range("L82").Value = 3212.5
$ws.Range("M82").Value = -111113339
$ws.Range("N82").Value = -3934.5
$ws.Range("H85").Value = 66669500
$ws.Range("I85").Value = 111113700
$ws.Range("J85").Value = 3212.5
$ws.Range("K85").Value = 111113700
$ws.Range("L85").Value = 3212.5
$ws.Range("M85").Value = -111112452
$ws.Range("N85").Value = -5708.5
$ws.Range("H122").Value = 3682.111
$ws.Range("I122").Value = 2610.2307
$ws.Range("J122").Value = 6469
$ws.Range("K122").Value = 7830.6921
$ws.Range("L122").Value = 19407
$ws.Range("M122").Value = -5380.6921
$ws.Range("N122").Value = -24307
$ws.Range("H126").Value = 6673.647
$ws.Range("I126").Value = 5332.9165
$ws.Range("K126").Value = 15998.7495
$ws.Range("M126").Value = -13528.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2194.7673
$ws.Range("I122").Value = 1476.2084
$ws.Range("J122").Value = 3102.4211
$ws.Range("K122").Value = 4428.6252
$ws.Range("L122").Value = 9307.263300000001
$ws.Range("M122").Value = -1978.6252
$ws.Range("N122").Value = -14207.2633
$ws.Range("H136").Value = 1811.48
$ws.Range("I136").Value = 1056.5714
$ws.Range("K136").Value = 3169.7142
$ws.Range("M136").Value = -619.7142000000003
